$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.668.53"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "'2.279.70"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'95.07"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "'267.10"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("D10").Value = "'44.20"
$ws.Range("E10").Value = "  -7.55%  "
$ws.Range("D11").Value = "'0.0934"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'7.71"
$ws.Range("E12").Value = "  -6.96%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "'2.622.96"
$ws.Range("D15").Value = "'15.16"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'0.845"
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("D17").Value = "'2.283.95"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "'43.588.18"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'71.99"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").Value = "'2.34"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "'234.63"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'11.39"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'2.49"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'38.87"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").Value = "'176.55"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'21.94"
$ws.Range("E32").Value = "  +4.23%  "
$ws.Range("D33").Value = "'0.0881"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "'0.0355"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("D38").Value = "'4.43"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("E39").Value = "  -9.12%  "
$ws.Range("E40").Value = "  +8.69%  "
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("E42").Value = "  +17.55%  "
$ws.Range("D43").Value = "'11.83"
$ws.Range("E43").Value = "  -5.16%  "
$ws.Range("D44").Value = "'62.32"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").Value = "'8.80"
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'98.30"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'1.52"
$ws.Range("E50").Value = "  +6.52%  "
$ws.Range("D51").Value = "'2.502.59"
$ws.Range("E51").Value = "  +1.94%  "
